$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.318.96"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3
$ws.Range("D3").Value = "1.877.57"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7202"
$ws.Range("E5").Value = "  +1.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.80"
$ws.Range("E6").Value = "  +0.46%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07993"
$ws.Range("E8").Value = "  +2.44%  "

# Row 9
$ws.Range("E9").Value = "  +1.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.96"
$ws.Range("E10").Value = "  -0.28%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08216"
$ws.Range("E11").Value = "  -2.13%  "

# Row 12
$ws.Range("D12").Value = "1.880.97"
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("E13").Value = "  +3.94%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.225"
$ws.Range("E14").Value = "  -0.24%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7117"
$ws.Range("E15").Value = "  +0.13%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.412"
$ws.Range("E16").Value = "  +5.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008488"
$ws.Range("E17").Value = "  +3.88%  "

# Row 18
$ws.Range("D18").Value = "29.322.54"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.39"
$ws.Range("E19").Value = "  +1.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.28"
$ws.Range("E20").Value = "  +0.40%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.754"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("E23").Value = "  +0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1597"
$ws.Range("E24").Value = "  +0.22%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.45"
$ws.Range("E25").Value = "  -0.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.037"
$ws.Range("E26").Value = "  +0.46%  "

# Row 27
$ws.Range("E27").Value = "  +0.29%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.502"
$ws.Range("E28").Value = "  -0.16%  "

# Row 29
$ws.Range("E29").Value = "  +0.28%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.304"
$ws.Range("E30").Value = "  +0.30%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.190"
$ws.Range("E31").Value = "  -8.25%  "

# Row 32
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("E33").Value = "  -0.59%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7602"
$ws.Range("E34").Value = "  +1.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("E36").Value = "  +0.65%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01871"
$ws.Range("E37").Value = "  +0.18%  "

# Row 38
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.276.44"
$ws.Range("E38").Value = "  +3.66%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.752"
$ws.Range("E39").Value = "  +0.91%  "

# Row 40
$ws.Range("E40").Value = "  -1.13%  "

# Row 41
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.80"
$ws.Range("E41").Value = "  +4.03%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9109"
$ws.Range("E42").Value = "  +2.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.11"
$ws.Range("E43").Value = "  +2.39%  "

# Row 44
$ws.Range("E44").Value = "  +7.28%  "

# Row 45
$ws.Range("E45").Value = "  +0.20%  "

# Row 46
$ws.Range("D46").Value = "2.025.18"
$ws.Range("E46").Value = "  +0.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5229"
$ws.Range("E47").Value = "  +0.64%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.794"
$ws.Range("E48").Value = "  +0.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.505"
$ws.Range("E49").Value = "  +1.00%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4337"
$ws.Range("E50").Value = "  +0.63%  "

# Row 51
$ws.Range("E51").Value = "  +0.29%  "
